# Regenerate merged AHB files
# - Rename header row labels: *_old -> *_FV2310, *_new -> *_FV2404
# - Freeze the header row (row 1)
# - Wrap the data range in an Excel Table (ListObject)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1, columns A:U) ---------------------
$ws.Range("A1").Value  = "Segmentname_FV2310"
$ws.Range("B1").Value  = "Segmentgruppe_FV2310"
$ws.Range("C1").Value  = "Segment_FV2310"
$ws.Range("D1").Value  = "Datenelement_FV2310"
$ws.Range("E1").Value  = "Segment ID_FV2310"
$ws.Range("F1").Value  = "Code_FV2310"
$ws.Range("G1").Value  = "Qualifier_FV2310"
$ws.Range("H1").Value  = "Beschreibung_FV2310"
$ws.Range("I1").Value  = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value  = "Bedingung_FV2310"
$ws.Range("K1").Value  = "diff"
$ws.Range("L1").Value  = "Segmentname_FV2404"
$ws.Range("M1").Value  = "Segmentgruppe_FV2404"
$ws.Range("N1").Value  = "Segment_FV2404"
$ws.Range("O1").Value  = "Datenelement_FV2404"
$ws.Range("P1").Value  = "Segment ID_FV2404"
$ws.Range("Q1").Value  = "Code_FV2404"
$ws.Range("R1").Value  = "Qualifier_FV2404"
$ws.Range("S1").Value  = "Beschreibung_FV2404"
$ws.Range("T1").Value  = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value  = "Bedingung_FV2404"

# --- 2. Freeze panes above row 2 (keeps header row visible) -------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the A1:U85 range into an Excel Table (ListObject) ----------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), $null, 1)
$tbl.Name = "Table1"
